$wb = $excel.ActiveWorkbook

# --- Update metadata on the "Metadata" worksheet ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/number-of-units-allowed"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Clear the Constraint(s) value for the top-level Extension row on the "Elements" worksheet ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
